# rating countries visited (plotly) + Beautify Folium map
#
# Inserts a new "Region_General" column (B) into the "Countries" sheet,
# fills it with a coarse Europe/Asia/Middle East grouping, fixes a data
# value (Germany's History rating), and refreshes the sort/conditional-
# formatting ranges that the column insert leaves stale.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Countries")

# --- 1. Insert a new column B ("Region_General") -------------------------
# Copy column A onto itself before inserting so the freed-up column keeps
# an exact (bit-for-bit) copy of the original stored width instead of a
# value re-derived (and re-rounded) through the ColumnWidth setter.
$ws.Columns("A:A").Copy()
$ws.Columns("A:A").Insert()

# --- 2. Header + values for the new column --------------------------------
$ws.Range("B1").Value = "Region_General"

$ws.Range("B2").Value = "Europe"
$ws.Range("B3").Value = "Europe"
$ws.Range("B4").Value = "Europe"
$ws.Range("B5").Value = "Europe"
$ws.Range("B6").Value = "Middle East"
$ws.Range("B7").Value = "Middle East"
$ws.Range("B8").Value = "Middle East"
$ws.Range("B9").Value = "Europe"
$ws.Range("B10").Value = "Asia"
$ws.Range("B11").Value = "Asia"
$ws.Range("B12").Value = "Asia"
$ws.Range("B13").Value = "Asia"
$ws.Range("B14").Value = "Europe"
$ws.Range("B15").Value = "Europe"
$ws.Range("B16").Value = "Europe"
$ws.Range("B17").Value = "Europe"

# --- 3. Data fix: Germany's History rating 5 -> 4 (now column J) ---------
$ws.Range("J16").Value = 4

# --- 4. Refresh the sort range (was A2:B17, now spans the extra column) --
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A17"))
$ws.Sort.SetRange($ws.Range("A1:C17"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- 5. Shift the color-scale conditional formatting one column right ----
# Re-point the existing rule instead of deleting/recreating it so its
# custom colors (F8696B / FCFCFF / 63BE7B) survive.
$cfRule = $ws.Range("G2:K17").FormatConditions.Item(1)
$cfRule.ModifyAppliesToRange($ws.Range("H2:L17"))

# --- 6. Selection matches the post-edit cursor position -------------------
$ws.Range("K10").Select()
